$wb = $excel.ActiveWorkbook

# Rename the two sheets (this will also update dependent formulas / defined names automatically)
$wsA = $wb.Worksheets.Item("scenarioA")
$wsB = $wb.Worksheets.Item("scenarioB")

$wsA.Name = "ScenarioA"
$wsB.Name = "ScenarioB"

# Fix the header labels in row 3 (M3/N3) on both sheets: "lat"/"lon" -> "Latitude"/"Longitude"
foreach ($ws in @($wsA, $wsB)) {
    # Copy M3's cell formatting onto N3 first (N3 should match M3's style)
    $ws.Range("M3").Copy()
    $ws.Range("N3").PasteSpecial(-4122)

    $ws.Range("M3").Value = "Latitude"
    $ws.Range("N3").Value = "Longitude"
}
